# Case and Fatality Demographics Data Updated
# Updates the "Fatalities by ..." breakdown sheets with refreshed counts.
# Percentage columns are formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook

# --- Fatalities by Age Group ---
$wsFatAge = $wb.Worksheets.Item("Fatalities by Age Group")
$wsFatAge.Range("B5").Value  = 274
$wsFatAge.Range("B6").Value  = 904
$wsFatAge.Range("B7").Value  = 2641
$wsFatAge.Range("B8").Value  = 5914
$wsFatAge.Range("B9").Value  = 4893
$wsFatAge.Range("B10").Value = 6275
$wsFatAge.Range("B11").Value = 6897
$wsFatAge.Range("B12").Value = 6791
$wsFatAge.Range("B13").Value = 16967

# --- Fatalities by Gender ---
$wsFatGender = $wb.Worksheets.Item("Fatalities by Gender")
$wsFatGender.Range("B2").Value = 21650
$wsFatGender.Range("B3").Value = 29959

# --- Fatalities by Race-Ethnicity ---
$wsFatRace = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$wsFatRace.Range("B2").Value = 1095
$wsFatRace.Range("B3").Value = 5274
$wsFatRace.Range("B4").Value = 23946
$wsFatRace.Range("B6").Value = 20988

# Make "Fatalities by Race-Ethnicity" the active/selected sheet, matching
# the refreshed report's tab focus.
$wsFatRace.Activate()
